# quarterly.xlsx update: roll the quarter window forward by one quarter
# (drop "Q4 1399/06", append new "Q2 1401/12") and refresh the read_price
# derived figures for every metric row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header / quarter-label rows (row 8 "هزینه های عمومی و اداری" table,
#    row 24 "تعداد پرسنل" table): shift each quarter label one column to
#    the left (E<-F<-G...<-M<-N) and put the brand-new quarter in N.
# ---------------------------------------------------------------------
$quarters = @(
    "فصل اول منتهی به 1399/09",
    "فصل دوم منتهی به 1399/12",
    "فصل سوم منتهی به 1400/03",
    "فصل چهارم منتهی به 1400/06",
    "فصل اول منتهی به 1400/09",
    "فصل دوم منتهی به 1400/12",
    "فصل سوم منتهی به 1401/03",
    "فصل چهارم منتهی به 1401/06",
    "فصل اول منتهی به 1401/09",
    "فصل دوم منتهی به 1401/12"
)

$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $quarters[$i] }
    $ws.Range("E$r`:N$r").Value = $arr
}

# ---------------------------------------------------------------------
# 2) Data rows: every metric shifts one quarter to the left, the newest
#    quarter (column N) receives the freshly read value.
# ---------------------------------------------------------------------
$rows = @{
    10 = @(28930,40928,82005,241417,229143,337374,243798,113190,101991,172230)
    11 = @(0,0,0,0,0,0,0,0,0,0)
    12 = @(0,0,0,0,0,0,0,0,0,0)
    13 = @(0,0,0,0,0,0,0,28766,0,0)
    14 = @(0,0,0,0,0,0,0,0,0,0)
    15 = @(313,428,282,710,23,32,794,-724,23,164)
    16 = @(2810,773,789,866,51,11,278,1608,814,1225)
    17 = @(17140,17555,23548,30789,20761,16682,29017,22621,19627,23329)
    18 = @(0,0,0,0,0,299,0,-299,0,0)
    19 = @(8480,14319,14123,33163,25837,45918,36691,19244,12419,-10571)
    20 = @(57673,74003,120747,306945,275815,400316,310578,184406,134874,186377)
    26 = @(43,43,47,42,48,48,35,41,41,42)
    27 = @(286,281,305,312,312,312,320,312,312,304)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("E$r`:N$r").Value = $arr
}

# ---------------------------------------------------------------------
# 3) Row-height touch-up that came along with the resave.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14.4
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(3).RowHeight = 14.4
$ws.Rows.Item(4).RowHeight = 14.4
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(7).RowHeight = 14.4
$ws.Rows.Item(8).RowHeight = 31.2
$ws.Rows.Item(9).RowHeight = 14.4
for ($r = 10; $r -le 23; $r++) { $ws.Rows.Item($r).RowHeight = 14.4 }
$ws.Rows.Item(24).RowHeight = 31.2
for ($r = 25; $r -le 31; $r++) { $ws.Rows.Item($r).RowHeight = 14.4 }
